$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the header text for columns I and J:
#   I1: "Jenis Instansi (Lokal/Nasional/Internasional)" -> "Skala Instansi (Lokal/Nasional/Internasional)"
#   J1: "Posisi" -> "Profesi"
$ws.Range("I1").Value = "Skala Instansi (Lokal/Nasional/Internasional)"
$ws.Range("J1").Value = "Profesi"
